$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    ,(73.0, 'Tuesday, Jan 10', '6:10 AM', 'FR1940', 'Bristol', '(BRS)', 'Ryanair ', 'B738', '(SP-RKT)', '6:13 AM', $null, '0 hours, 3 minutes', $null)
    ,(74.0, 'Tuesday, Jan 10', '6:10 AM', 'FR2008', 'London', '(STN)', 'Ryanair ', 'B738', '(SP-RKL)', '6:14 AM', $null, '0 hours, 4 minutes', $null)
    ,(75.0, 'Tuesday, Jan 10', '6:20 AM', 'FR1963', 'Milan', '(BGY)', 'Ryanair ', 'B738', '(SP-RSV)', '6:23 AM', $null, '0 hours, 3 minutes', $null)
    ,(76.0, 'Tuesday, Jan 10', '7:00 AM', 'FR1888', 'Paris', '(BVA)', 'Ryanair ', 'B38M', '(SP-RZO)', '7:00 AM', $null, '0 hours, 0 minutes', $null)
    ,(77.0, 'Tuesday, Jan 10', '7:05 AM', 'FR4177', 'Pafos', '(PFO)', 'Ryanair ', 'B38M', '(SP-RZI)', '7:04 AM', $null, '0 hours, -1 minutes', $null)
    ,(78.0, 'Tuesday, Jan 10', '7:10 AM', 'FR1115', 'Rome', '(CIA)', 'Ryanair ', 'B738', '(SP-RKP)', '7:09 AM', $null, '0 hours, -1 minutes', $null)
    ,(79.0, 'Tuesday, Jan 10', '9:30 AM', 'FR1922', 'Budapest', '(BUD)', 'Ryanair ', 'B738', '(SP-RSC)', '9:30 AM', $null, '0 hours, 0 minutes', $null)
    ,(80.0, 'Tuesday, Jan 10', '10:30 AM', 'FR8012', 'Dublin', '(DUB)', 'Ryanair ', 'B38M', '(EI-HGX)', '10:48 AM', $null, '0 hours, 18 minutes', $null)
    ,(81.0, 'Tuesday, Jan 10', '11:35 AM', 'FR2263', 'Lisbon', '(LIS)', 'Ryanair ', 'B738', '(SP-RSP)', '11:46 AM', $null, '0 hours, 11 minutes', $null)
    ,(82.0, 'Tuesday, Jan 10', '12:30 PM', 'FR4060', 'Malaga', '(AGP)', 'Buzz ', 'B38M', '(SP-RZG)', '12:32 PM', $null, '0 hours, 2 minutes', $null)
    ,(83.0, 'Tuesday, Jan 10', '12:35 PM', 'FR2670', 'London', '(STN)', 'Ryanair ', 'B738', '(SP-RKT)', '1:32 PM', $null, '0 hours, 57 minutes', $null)
    ,(84.0, 'Tuesday, Jan 10', '2:55 PM', 'FR1938', 'Gothenburg', '(GOT)', 'Ryanair ', 'B38M', '(SP-RZI)', '3:00 PM', $null, '0 hours, 5 minutes', $null)
    ,(85.0, 'Tuesday, Jan 10', '6:25 PM', 'FR1968', 'Madrid', '(MAD)', 'Ryanair ', 'B738', '(EI-ENL)', '7:49 PM', $null, '1 hours, 24 minutes', $null)
    ,(86.0, 'Tuesday, Jan 10', '6:45 PM', 'FR1022', 'London', '(STN)', 'Ryanair ', 'B38M', '(SP-RZO)', '6:51 PM', $null, '0 hours, 6 minutes', $null)
    ,(87.0, 'Tuesday, Jan 10', '6:50 PM', 'FR1574', 'Vienna', '(VIE)', 'Ryanair ', 'B38M', '(SP-RZI)', '7:07 PM', $null, '0 hours, 17 minutes', $null)
    ,(88.0, 'Tuesday, Jan 10', '6:55 PM', 'FR1904', 'Milan', '(BGY)', 'Ryanair ', 'B738', '(SP-RKP)', '7:23 PM', $null, '0 hours, 28 minutes', $null)
    ,(89.0, 'Tuesday, Jan 10', '7:50 PM', 'UNKNOWN', 'Belgrade', '(BEG)', 'AMC Aviation ', 'LJ60', '(SP-CEZ)', '7:46 PM', $null, '0 hours, -4 minutes', $null)
    ,(90.0, 'Tuesday, Jan 10', '10:00 PM', 'UNKNOWN', 'Belgrade', '(BEG)', 'AMC Aviation ', 'LJ60', '(SP-CEZ)', '10:38 PM', $null, '0 hours, 38 minutes', $null)
)

$startRow = 74
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le 13; $c++) {
        $val = $rowData[$c - 1]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}
